$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: convert Meta/Meta.AC/Venda/Venda.AC/Pecas.AC/Sobras/P from text to real numbers ---
$ws.Cells.Item(3, 2).Value = 2525       # B3
$ws.Cells.Item(3, 3).Value = 7535       # C3
$ws.Cells.Item(3, 4).Value = 2412.4     # D3
$ws.Cells.Item(3, 5).Value = 6626.4     # E3
$ws.Cells.Item(3, 6).Value = 54         # F3
$ws.Cells.Item(3, 7).Value = 908.6      # G3
$ws.Cells.Item(3, 8).Value = 87.94      # H3

# --- Row 4: new data row, numeric values ---
$ws.Cells.Item(4, 1).Value = "31/35/5000"
$ws.Cells.Item(4, 2).Value = 5000
$ws.Cells.Item(4, 3).Value = 12535
$ws.Cells.Item(4, 4).Value = 5000
$ws.Cells.Item(4, 5).Value = 11626.4
$ws.Cells.Item(4, 6).Value = 74
$ws.Cells.Item(4, 7).Value = 908.6
$ws.Cells.Item(4, 8).Value = 92.75

# --- Row 5: new data row, numeric values ---
$ws.Cells.Item(5, 1).Value = "31/06/2000"
$ws.Cells.Item(5, 2).Value = 2000
$ws.Cells.Item(5, 3).Value = 14535
$ws.Cells.Item(5, 4).Value = 2000
$ws.Cells.Item(5, 5).Value = 13626.4
$ws.Cells.Item(5, 6).Value = 94
$ws.Cells.Item(5, 7).Value = 908.6
$ws.Cells.Item(5, 8).Value = 93.75

# --- Row 6: new data row, values kept as formatted text (e.g. "8000.00") ---
$ws.Cells.Item(6, 1).Value = "20/08/8000"

$textCells = @(
    @{ Col = 2; Text = "8000.00" },
    @{ Col = 3; Text = "22535.00" },
    @{ Col = 4; Text = "8000.00" },
    @{ Col = 5; Text = "21626.40" },
    @{ Col = 6; Text = "144.0" },
    @{ Col = 7; Text = "908.60" },
    @{ Col = 8; Text = "95.97" }
)

foreach ($tc in $textCells) {
    $cell = $ws.Cells.Item(6, $tc.Col)
    $cell.Value = "'" + $tc.Text
    $cell.Style = "Normal"
}
